$d = $word.ActiveDocument

# 1. Add "(Kick-off)" after "P4P SubscribeMe "
$r = $d.Content
$r.Find.Execute("P4P SubscribeMe ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter("(Kick-off)")

# 2. Remove the word "verder " from the long paragraph
$d.Content.Find.Execute("projectomschrijving verder toegelicht", $true, $false, $false, $false, $false,
                         $true, 1, $false, "projectomschrijving toegelicht", 2)
